$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-03 Saturday" "2024-08-04 Sunday"

Replace-Text "11÷5=" "43÷8="
Replace-Text "28÷3=" "40÷8="
Replace-Text "68÷8=" "95÷5="
Replace-Text "21÷2=" "69÷9="
Replace-Text "98÷4=" "46÷9="
Replace-Text "64÷7=" "76÷5="
Replace-Text "88÷8=" "79÷6="
Replace-Text "11÷7=" "29÷7="
Replace-Text "82÷3=" "14÷7="
Replace-Text "28÷8=" "66÷4="
Replace-Text "93÷2=" "41÷2="
Replace-Text "22÷2=" "62÷8="
Replace-Text "20÷6=" "50÷3="
Replace-Text "17÷9=" "66÷3="
Replace-Text "49÷4=" "69÷6="
Replace-Text "48÷3=" "46÷9="
Replace-Text "67÷6=" "78÷4="
Replace-Text "61÷2=" "90÷5="
Replace-Text "19÷5=" "21÷3="
Replace-Text "11÷3=" "14÷3="
Replace-Text "98÷9=" "85÷2="
Replace-Text "60÷6=" "19÷7="
Replace-Text "47÷6=" "78÷6="
Replace-Text "99÷5=" "18÷4="
Replace-Text "23÷7=" "76÷6="
